# Insert a new weekly record as row 192 in the "Acelga" price sheet.
# This shifts the previous rows 192:202 down to 193:203 and fills the
# newly opened row 192 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 192, pushing existing rows down.
$ws.Rows.Item(192).Insert()

# Populate the new row 192 with the new record (same template values as
# the other rows for the fixed columns, new values for the variable ones).
$ws.Cells.Item(192, 1).Value = 5
$ws.Cells.Item(192, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(192, 3).Value = "Maule"
$ws.Cells.Item(192, 4).Value = 44516
$ws.Cells.Item(192, 5).Value = 7
$ws.Cells.Item(192, 6).Value = 100112009
$ws.Cells.Item(192, 7).Value = "Acelga"
$ws.Cells.Item(192, 8).Value = "Sin especificar"
$ws.Cells.Item(192, 9).Value = "Primera"
$ws.Cells.Item(192, 10).Value = 500
$ws.Cells.Item(192, 11).Value = 2000
$ws.Cells.Item(192, 12).Value = 2000
$ws.Cells.Item(192, 13).Value = 2000
$ws.Cells.Item(192, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(192, 15).Value = "Región del Maule"
$ws.Cells.Item(192, 16).Value = 500
$ws.Cells.Item(192, 17).Value = 4
$ws.Cells.Item(192, 18).Value = "Hortaliza"
